$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Insert a new row at 68 (pushes old rows 68-74 down to 69-75). xlShiftDown
# with CopyOrigin=xlFormatFromLeftOrAbove copies row 67's cell formatting
# onto the fresh row (A/B get styles 2/3, matching what row 67 already has).
$ws.Rows.Item(68).Insert(-4121, 0) | Out-Null
$ws.Rows.Item(68).RowHeight = 15.75

$ws.Range("A68").Value = "test_instructions_4"
$ws.Range("B68").Value = "For best results, do not sing along with your instrument. Just play!"
$ws.Range("C68").Value = "For best results, do not sing along with your instrument. Just play!"

# C68 should carry the same (fill-shaded) style as B68 rather than the
# plain style it inherited from row 67's column C.
$ws.Range("B68").Copy() | Out-Null
$ws.Range("C68").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Reflect where the author ended up after making the edit.
$ws.Range("A68").Select() | Out-Null
